$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1.xml) -- table "Overview"
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 5: 1e4ff530-9bc6-4869-9acc-192cd47c1999.md
$wsOverview.Range("A5").Value = "1e4ff530-9bc6-4869-9acc-192cd47c1999.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/1e4ff530-9bc6-4869-9acc-192cd47c1999.md", "", "", "e2e\1e4ff530-9bc6-4869-9acc-192cd47c1999.md")
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-06 12:17:14"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 6: 27087249-2b9e-4d49-afd7-ca98a7bcdbac.png
$wsOverview.Range("A6").Value = "27087249-2b9e-4d49-afd7-ca98a7bcdbac.png"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/27087249-2b9e-4d49-afd7-ca98a7bcdbac.png", "", "", "e2e\27087249-2b9e-4d49-afd7-ca98a7bcdbac.png")
$wsOverview.Range("C6").Value = ".png"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-09-06 12:17:14"
$wsOverview.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 7: dca2d15c-9522-4382-9bfd-58768820b51e.png
$wsOverview.Range("A7").Value = "dca2d15c-9522-4382-9bfd-58768820b51e.png"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/dca2d15c-9522-4382-9bfd-58768820b51e.png", "", "", "e2e\dca2d15c-9522-4382-9bfd-58768820b51e.png")
$wsOverview.Range("C7").Value = ".png"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-09-06 12:17:14"
$wsOverview.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2.xml) -- table "zh-cn"
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, I=Latest Target File,
#          J=Latest Handback File, K=Latest Handback DateTime,
#          L=Reference Tokens, M=To be localized, N=Dependency From,
#          O=Has metadata, P=Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

# Row 5: 1e4ff530-9bc6-4869-9acc-192cd47c1999.md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/1e4ff530-9bc6-4869-9acc-192cd47c1999.md", "", "", "1e4ff530-9bc6-4869-9acc-192cd47c1999.md")
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "1e4ff530-9bc6-4869-9acc-192cd47c1999.6a79650514ebfaefe82dad92337e8c7a5eb934f2.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-06 12:16:57"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("O5").Value = "False"

# Row 6: 27087249-2b9e-4d49-afd7-ca98a7bcdbac.png
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/27087249-2b9e-4d49-afd7-ca98a7bcdbac.png", "", "", "27087249-2b9e-4d49-afd7-ca98a7bcdbac.png")
$wsZhCn.Range("B6").Value = ".png"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "ab5dc317ad9f8991cd7530d4b32775162e0331c5.png"
$wsZhCn.Range("H6").Value = "2016-09-06 12:16:57"
$wsZhCn.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M6").Value = "True(Dependency)"
$wsZhCn.Range("N6").Value = "e2e\1e4ff530-9bc6-4869-9acc-192cd47c1999.md"
$wsZhCn.Range("O6").Value = "False"

# Row 7: dca2d15c-9522-4382-9bfd-58768820b51e.png
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/dca2d15c-9522-4382-9bfd-58768820b51e.png", "", "", "dca2d15c-9522-4382-9bfd-58768820b51e.png")
$wsZhCn.Range("B7").Value = ".png"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "7989491fc9144ac5d33f31357ed103fd91bdb572.png"
$wsZhCn.Range("H7").Value = "2016-09-06 12:16:57"
$wsZhCn.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M7").Value = "True(Dependency)"
$wsZhCn.Range("N7").Value = "e2e\1e4ff530-9bc6-4869-9acc-192cd47c1999.md"
$wsZhCn.Range("O7").Value = "False"

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3.xml) -- table "de-de"
# Same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

# Row 5: 1e4ff530-9bc6-4869-9acc-192cd47c1999.md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/1e4ff530-9bc6-4869-9acc-192cd47c1999.md", "", "", "1e4ff530-9bc6-4869-9acc-192cd47c1999.md")
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "1e4ff530-9bc6-4869-9acc-192cd47c1999.6a79650514ebfaefe82dad92337e8c7a5eb934f2.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-06 12:17:14"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("O5").Value = "False"

# Row 6: 27087249-2b9e-4d49-afd7-ca98a7bcdbac.png
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/27087249-2b9e-4d49-afd7-ca98a7bcdbac.png", "", "", "27087249-2b9e-4d49-afd7-ca98a7bcdbac.png")
$wsDeDe.Range("B6").Value = ".png"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "ab5dc317ad9f8991cd7530d4b32775162e0331c5.png"
$wsDeDe.Range("H6").Value = "2016-09-06 12:17:14"
$wsDeDe.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M6").Value = "True(Dependency)"
$wsDeDe.Range("N6").Value = "e2e\1e4ff530-9bc6-4869-9acc-192cd47c1999.md"
$wsDeDe.Range("O6").Value = "False"

# Row 7: dca2d15c-9522-4382-9bfd-58768820b51e.png
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e472a28ce674153e17ce51e5369a5d23059096d3/e2e/dca2d15c-9522-4382-9bfd-58768820b51e.png", "", "", "dca2d15c-9522-4382-9bfd-58768820b51e.png")
$wsDeDe.Range("B7").Value = ".png"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "7989491fc9144ac5d33f31357ed103fd91bdb572.png"
$wsDeDe.Range("H7").Value = "2016-09-06 12:17:14"
$wsDeDe.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M7").Value = "True(Dependency)"
$wsDeDe.Range("N7").Value = "e2e\1e4ff530-9bc6-4869-9acc-192cd47c1999.md"
$wsDeDe.Range("O7").Value = "False"
